$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: 153. Find Minimum in Rotated Sorted Array
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null

$ws.Range("A8").Value = "153. Find Minimum in Rotated Sorted Array"
$ws.Range("B8").Value = "Medium"
$ws.Range("C8").Value = "Binary Search"
$ws.Range("D8").Value = "We consider where the pivot is in left and right portions, and need a basic check for an already sorted portion."
$ws.Range("E8").Value = "https://leetcode.com/problems/find-minimum-in-rotated-sorted-array/solutions/158940/beat-100-very-simple-python-very-detailed-explanation/ "
$ws.Hyperlinks.Add($ws.Range("E8"), "https://leetcode.com/problems/find-minimum-in-rotated-sorted-array/solutions/158940/beat-100-very-simple-python-very-detailed-explanation/") | Out-Null

# Row 9: 22. Generate Parentheses
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null

$ws.Range("A9").Value = "22. Generate Parentheses"
$ws.Range("B9").Value = "Medium"
$ws.Range("C9").Value = "Stack"
$ws.Range("D9").Value = "Custom class and stack, track string, # of open and closed parentheses per step. Use DFS loop."
$ws.Range("E9").Value = "https://leetcode.com/problems/generate-parentheses/solutions/10391/java-solution-using-stack/ "
$ws.Hyperlinks.Add($ws.Range("E9"), "https://leetcode.com/problems/generate-parentheses/solutions/10391/java-solution-using-stack/") | Out-Null

# Resize table to include new rows
$ws.ListObjects.Item(1).Resize($ws.Range("A1:E9"))

# Adjust column A width to fit new (longer) content (bestFit-equivalent)
$ws.Columns.Item(1).ColumnWidth = 39.0

# Update selection to match target state
$ws.Range("D11").Select() | Out-Null

$wb.Save()
